$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cell = $ws.Range("A4")
$cell.Value = 43252
$cell.NumberFormat = "mm-dd-yy"
$cell.Font.Name = "Calibri"
$cell.Font.Size = 11
$cell.Borders.Color = 0
$cell.Borders.LineStyle = 1
$cell.Borders.Weight = 2
$cell.Interior.Color = 15986394

$ws.Range("A4").Select()
$ws.Columns.Item(1).ColumnWidth = 10.71364
